$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 26858
$ws1.Range("F4").Value = 594
$ws1.Range("G5").Value = 68
$ws1.Range("F6").Value = 621
$ws1.Range("F7").Value = 179
$ws1.Range("F10").Value = 366
$ws1.Range("F11").Value = 453
$ws1.Range("F12").Value = 192
$ws1.Range("F13").Value = 51
$ws1.Range("F14").Value = 305
$ws1.Range("F15").Value = 78
$ws1.Range("F16").Value = 445
$ws1.Range("F18").Value = 1570
$ws1.Range("F19").Value = 221
$ws1.Range("F20").Value = 55

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 4513
$ws2.Range("F3").Value = 238
$ws2.Range("F6").Value = 206
$ws2.Range("F7").Value = 206
$ws2.Range("F10").Value = 442
$ws2.Range("F15").Value = 66

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5115
$ws3.Range("F3").Value = 246

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5115
$ws4.Range("F4").Value = 246
$ws4.Range("F5").Value = 26858
$ws4.Range("F6").Value = 594
$ws4.Range("F7").Value = 4513
$ws4.Range("G8").Value = 68
$ws4.Range("F9").Value = 238
$ws4.Range("F10").Value = 621
$ws4.Range("F13").Value = 179
$ws4.Range("F14").Value = 206
$ws4.Range("F15").Value = 206
$ws4.Range("F18").Value = 442
$ws4.Range("F22").Value = 366
$ws4.Range("F23").Value = 453
$ws4.Range("F24").Value = 192
$ws4.Range("F25").Value = 51
$ws4.Range("F27").Value = 305
$ws4.Range("F28").Value = 78
$ws4.Range("F31").Value = 445
$ws4.Range("F33").Value = 66
$ws4.Range("F34").Value = 1570
$ws4.Range("F35").Value = 221
$ws4.Range("F37").Value = 55
